$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the LaunchTime text values for balloonA/B/C (date changed from 27 Feb 2018 to 18 Mar 2018).
# A leading apostrophe forces Excel to store the value as literal text (quote-prefixed),
# matching the original cell formatting; Excel strips the leading apostrophe marker itself.
$ws.Range("B2").Value = "'18 Mar 2018 17:00:00.000'"
$ws.Range("B3").Value = "'18 Mar 2018 18:00:00.000'"
$ws.Range("B4").Value = "'18 Mar 2018 19:00:00.000'"

# Move the active selection from F5 to B5
$ws.Range("B5").Select()
